# Auto-generated Excel COM-interop edit script
# Applies scheduled market-price / profit recalculation updates to the
# Leve profit tables across all eight crafting-job worksheets.
#
# For each touched row this sets the new currentAveragePrice* / LevePrice* /
# LeveProfit* values (columns H-N). A few rows also gain or lose a trailing
# LeveProfitHQ (column N) or LeveProfitNQ (column M) cell entirely, matching
# the source diff exactly (ClearContents to drop a cell, plain assignment to
# add one).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4136.1
$ws.Range("I62").Value = 3725
$ws.Range("J62").Value = 4238.875
$ws.Range("K62").Value = 3725
$ws.Range("L62").Value = 4238.875
$ws.Range("M62").Value = -3101
$ws.Range("N62").Value = -5486.875
$ws.Range("H65").Value = 4136.1
$ws.Range("I65").Value = 3725
$ws.Range("J65").Value = 4238.875
$ws.Range("K65").Value = 18625
$ws.Range("L65").Value = 21194.375
$ws.Range("M65").Value = -15505
$ws.Range("N65").Value = -27434.375
$ws.Range("H70").Value = 1489.8
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1489.8
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4469.4
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -5009.4
$ws.Range("H73").Value = 1489.8
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1489.8
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4469.4
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -6341.4
$ws.Range("H74").Value = 4756.6875
$ws.Range("I74").Value = 3158.1428
$ws.Range("K74").Value = 3158.1428
$ws.Range("M74").Value = -2222.1428
$ws.Range("H77").Value = 4756.6875
$ws.Range("I77").Value = 3158.1428
$ws.Range("K77").Value = 15790.714
$ws.Range("M77").Value = -11110.714
$ws.Range("H100").Value = 2685.1538
$ws.Range("I100").Value = 1375.375
$ws.Range("K100").Value = 1375.375
$ws.Range("M100").Value = -834.375
$ws.Range("H131").Value = 3034.2727
$ws.Range("I131").Value = 2537.7
$ws.Range("K131").Value = 7613.099999999999
$ws.Range("M131").Value = -2573.099999999999
$ws.Range("H132").Value = 2365.575
$ws.Range("I132").Value = 1762.1212
$ws.Range("J132").Value = 5210.4287
$ws.Range("K132").Value = 5286.363600000001
$ws.Range("L132").Value = 15631.2861
$ws.Range("M132").Value = -2756.363600000001
$ws.Range("N132").Value = -20691.2861
$ws.Range("H138").Value = 11499634
$ws.Range("I138").Value = 1718.875
$ws.Range("J138").Value = 15879792
$ws.Range("K138").Value = 5156.625
$ws.Range("L138").Value = 47639376
$ws.Range("M138").Value = -16.625
$ws.Range("N138").Value = -47649656
$ws.Range("H141").Value = 2305
$ws.Range("J141").Value = 2553
$ws.Range("L141").Value = 7659
$ws.Range("N141").Value = -18019

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4446.778
$ws.Range("I32").Value = 4184.2046
$ws.Range("K32").Value = 4184.2046
$ws.Range("M32").Value = -3897.2046
$ws.Range("H45").Value = 1686.3182
$ws.Range("I45").Value = 1527.421
$ws.Range("K45").Value = 1527.421
$ws.Range("M45").Value = -1150.421
$ws.Range("H76").Value = 112729
$ws.Range("J76").Value = 112729
$ws.Range("L76").Value = 112729
$ws.Range("N76").Value = -113405
$ws.Range("H79").Value = 112729
$ws.Range("J79").Value = 112729
$ws.Range("L79").Value = 112729
$ws.Range("N79").Value = -115069
$ws.Range("H92").Value = 66663.336
$ws.Range("J92").Value = 66663.336
$ws.Range("L92").Value = 66663.336
$ws.Range("N92").Value = -71655.336
$ws.Range("H97").Value = 1801.3182
$ws.Range("I97").Value = 1799.6154
$ws.Range("J97").Value = 1803.7778
$ws.Range("K97").Value = 1799.6154
$ws.Range("L97").Value = 1803.7778
$ws.Range("M97").Value = -1303.6154
$ws.Range("N97").Value = -2795.7778
$ws.Range("H104").Value = 6750
$ws.Range("J104").Value = 6750
$ws.Range("L104").Value = 6750
$ws.Range("N104").Value = -13738
$ws.Range("H110").Value = 13643.621
$ws.Range("I110").Value = 15465.625
$ws.Range("K110").Value = 15465.625
$ws.Range("M110").Value = -13420.625
$ws.Range("H132").Value = 35767240
$ws.Range("I132").Value = 3216.2
$ws.Range("K132").Value = 9648.599999999999
$ws.Range("M132").Value = -7118.599999999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2150.8333
$ws.Range("I20").Value = 2215.5
$ws.Range("K20").Value = 2215.5
$ws.Range("M20").Value = -1968.5
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H82").Value = 41901.668
$ws.Range("I82").Value = 17877.5
$ws.Range("J82").Value = 89950
$ws.Range("K82").Value = 17877.5
$ws.Range("L82").Value = 89950
$ws.Range("M82").Value = -17494.5
$ws.Range("N82").Value = -90716
$ws.Range("H85").Value = 41901.668
$ws.Range("I85").Value = 17877.5
$ws.Range("J85").Value = 89950
$ws.Range("K85").Value = 17877.5
$ws.Range("L85").Value = 89950
$ws.Range("M85").Value = -16551.5
$ws.Range("N85").Value = -92602
$ws.Range("H86").Value = 13455.083
$ws.Range("I86").Value = 6542.3335
$ws.Range("J86").Value = 34193.332
$ws.Range("K86").Value = 6542.3335
$ws.Range("L86").Value = 34193.332
$ws.Range("M86").Value = -5419.3335
$ws.Range("N86").Value = -36439.332
$ws.Range("H89").Value = 13455.083
$ws.Range("I89").Value = 6542.3335
$ws.Range("J89").Value = 34193.332
$ws.Range("K89").Value = 32711.6675
$ws.Range("L89").Value = 170966.66
$ws.Range("M89").Value = -27095.6675
$ws.Range("N89").Value = -182198.66
$ws.Range("H94").Value = 1113.2941
$ws.Range("J94").Value = 2030.75
$ws.Range("L94").Value = 2030.75
$ws.Range("N94").Value = -2932.75
$ws.Range("H105").Value = 11337.182
$ws.Range("I105").Value = 14276.125
$ws.Range("K105").Value = 14276.125
$ws.Range("M105").Value = -12529.125
$ws.Range("H134").Value = 2837.625
$ws.Range("I134").Value = 2432.32
$ws.Range("J134").Value = 4285.143
$ws.Range("K134").Value = 7296.960000000001
$ws.Range("L134").Value = 12855.429
$ws.Range("M134").Value = -4761.960000000001
$ws.Range("N134").Value = -17925.429

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 529.25
$ws.Range("I16").Value = 464.53333
$ws.Range("K16").Value = 464.53333
$ws.Range("M16").Value = -177.53333
$ws.Range("H22").Value = 18035.334
$ws.Range("I22").Value = 20542.4
$ws.Range("K22").Value = 20542.4
$ws.Range("M22").Value = -20192.4
$ws.Range("H113").Value = 529.25
$ws.Range("I113").Value = 464.53333
$ws.Range("K113").Value = 464.53333
$ws.Range("M113").Value = 1705.46667
$ws.Range("H122").Value = 2239.2307
$ws.Range("I122").Value = 1748.1818
$ws.Range("K122").Value = 5244.5454
$ws.Range("M122").Value = -2794.5454
$ws.Range("H132").Value = 5108.1904
$ws.Range("I132").Value = 4237.3887
$ws.Range("K132").Value = 12712.1661
$ws.Range("M132").Value = -10182.1661

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 594.5454999999999
$ws.Range("I39").Value = 304.05
$ws.Range("J39").Value = 3499.5
$ws.Range("K39").Value = 912.1500000000001
$ws.Range("L39").Value = 10498.5
$ws.Range("M39").Value = -618.1500000000001
$ws.Range("N39").Value = -11086.5
$ws.Range("H113").Value = 2942.8235
$ws.Range("I113").Value = 1432
$ws.Range("K113").Value = 4296
$ws.Range("M113").Value = -2126
$ws.Range("H131").Value = 22625.943
$ws.Range("J131").Value = 4378.909
$ws.Range("L131").Value = 13136.727
$ws.Range("N131").Value = -23216.727
$ws.Range("H140").Value = 2097.0833
$ws.Range("I140").Value = 1909.9
$ws.Range("K140").Value = 5729.700000000001
$ws.Range("M140").Value = -549.7000000000007

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 50000.5
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 50000.5
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H102").Value = 2230.2104
$ws.Range("I102").Value = 619.3333
$ws.Range("K102").Value = 619.3333
$ws.Range("M102").Value = 1002.6667
$ws.Range("H113").Value = 2785.3447
$ws.Range("I113").Value = 1448.8334
$ws.Range("K113").Value = 1448.8334
$ws.Range("M113").Value = 721.1666
$ws.Range("H126").Value = 8668.799999999999
$ws.Range("I126").Value = 10337.4
$ws.Range("J126").Value = 7000.2
$ws.Range("K126").Value = 31012.2
$ws.Range("L126").Value = 21000.6
$ws.Range("M126").Value = -28542.2
$ws.Range("N126").Value = -25940.6

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2531.8076
$ws.Range("I40").Value = 2473.6191
$ws.Range("J40").Value = 2776.2
$ws.Range("K40").Value = 2473.6191
$ws.Range("L40").Value = 2776.2
$ws.Range("M40").Value = -2337.6191
$ws.Range("N40").Value = -3048.2
$ws.Range("H96").Value = 51862.332
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 51862.332
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 51862.332
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -57354.332

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7897.5
$ws.Range("J62").Value = 7897.5
$ws.Range("L62").Value = 7897.5
$ws.Range("N62").Value = -9145.5
$ws.Range("H65").Value = 7897.5
$ws.Range("J65").Value = 7897.5
$ws.Range("L65").Value = 39487.5
$ws.Range("N65").Value = -45727.5
$ws.Range("H99").Value = 59999
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 59999
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 59999
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -65989
$ws.Range("H100").Value = 166668720
$ws.Range("H126").Value = 7201.4165
$ws.Range("I126").Value = 7379.6665
$ws.Range("K126").Value = 22138.9995
$ws.Range("M126").Value = -19668.9995

